# Auto-generated script to apply 2023-06-11 daily crime data update
# across the Citywide Totals, By Neighborhood, and individual neighborhood sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("J2").Value = 3176
$ws.Range("J3").Value = 3325
$ws.Range("H4").Value = 1696
$ws.Range("J4").Value = 741
$ws.Range("J5").Value = 260
$ws.Range("J6").Value = 3917
$ws.Range("H7").Value = 26006
$ws.Range("J7").Value = 11419

$ws = $wb.Worksheets.Item('Uptown')
$ws.Range("J6").Value = 38
$ws.Range("J7").Value = 131

$ws = $wb.Worksheets.Item('West Ridge')
$ws.Range("J6").Value = 38
$ws.Range("J7").Value = 126

$ws = $wb.Worksheets.Item('Bridgeport')
$ws.Range("J3").Value = 10
$ws.Range("J7").Value = 48

$ws = $wb.Worksheets.Item('Fuller Park')
$ws.Range("J6").Value = 12
$ws.Range("J7").Value = 48

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range("J3").Value = 129
$ws.Range("J7").Value = 365

$ws = $wb.Worksheets.Item('Woodlawn')
$ws.Range("J2").Value = 52
$ws.Range("J6").Value = 46
$ws.Range("J7").Value = 167

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range("J5").Value = 10
$ws.Range("J6").Value = 115
$ws.Range("J7").Value = 406

$ws = $wb.Worksheets.Item('New City')
$ws.Range("J2").Value = 87
$ws.Range("J3").Value = 90
$ws.Range("J6").Value = 108
$ws.Range("J7").Value = 305

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("J7").Value = 346
$ws.Range("J8").Value = 732
$ws.Range("J14").Value = 48
$ws.Range("J19").Value = 359
$ws.Range("J20").Value = 240
$ws.Range("J22").Value = 25
$ws.Range("J23").Value = 117
$ws.Range("J25").Value = 65
$ws.Range("J27").Value = 69
$ws.Range("J29").Value = 655
$ws.Range("J30").Value = 48
$ws.Range("J33").Value = 495
$ws.Range("J34").Value = 58
$ws.Range("J37").Value = 365
$ws.Range("J42").Value = 457
$ws.Range("J43").Value = 105
$ws.Range("J48").Value = 112
$ws.Range("J50").Value = 64
$ws.Range("J51").Value = 154
$ws.Range("J52").Value = 310
$ws.Range("J53").Value = 111
$ws.Range("J54").Value = 216
$ws.Range("J57").Value = 52
$ws.Range("H63").Value = 248
$ws.Range("J63").Value = 49
$ws.Range("J65").Value = 305
$ws.Range("J66").Value = 27
$ws.Range("J67").Value = 406
$ws.Range("J71").Value = 42
$ws.Range("J72").Value = 46
$ws.Range("J78").Value = 153
$ws.Range("J79").Value = 340
$ws.Range("J83").Value = 260
$ws.Range("J85").Value = 520
$ws.Range("J89").Value = 131
$ws.Range("J91").Value = 128
$ws.Range("J94").Value = 103
$ws.Range("J95").Value = 180
$ws.Range("J96").Value = 126
$ws.Range("J97").Value = 68
$ws.Range("J99").Value = 167
$ws.Range("H101").Value = 26006
$ws.Range("J101").Value = 11419

$ws = $wb.Worksheets.Item('South Chicago')
$ws.Range("J2").Value = 75
$ws.Range("J7").Value = 260

$ws = $wb.Worksheets.Item('West Pullman')
$ws.Range("J2").Value = 68
$ws.Range("J3").Value = 54
$ws.Range("J6").Value = 46
$ws.Range("J7").Value = 180

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range("J2").Value = 132
$ws.Range("J3").Value = 159
$ws.Range("J7").Value = 495

$ws = $wb.Worksheets.Item('Loop')
$ws.Range("J2").Value = 55
$ws.Range("J7").Value = 216

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range("J3").Value = 226
$ws.Range("J4").Value = 39
$ws.Range("J7").Value = 655

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range("J2").Value = 84
$ws.Range("J3").Value = 102
$ws.Range("J7").Value = 359

$ws = $wb.Worksheets.Item('Lake View')
$ws.Range("J6").Value = 54
$ws.Range("J7").Value = 112

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range("J2").Value = 128
$ws.Range("J3").Value = 196
$ws.Range("J7").Value = 520

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range("J3").Value = 102
$ws.Range("J5").Value = 12
$ws.Range("J7").Value = 457

$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Range("J2").Value = 36
$ws.Range("J7").Value = 153

$ws = $wb.Worksheets.Item('Douglas')
$ws.Range("J3").Value = 41
$ws.Range("J7").Value = 117

$ws = $wb.Worksheets.Item('Washington Park')
$ws.Range("J3").Value = 60
$ws.Range("J7").Value = 128

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range("J3").Value = 123
$ws.Range("J6").Value = 94
$ws.Range("J7").Value = 340

$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Range("J5").Value = 4
$ws.Range("J7").Value = 240

$ws = $wb.Worksheets.Item('Little Village')
$ws.Range("J2").Value = 70
$ws.Range("J3").Value = 89
$ws.Range("J7").Value = 310

$ws = $wb.Worksheets.Item('Garfield Ridge')
$ws.Range("J6").Value = 22
$ws.Range("J7").Value = 58

$ws = $wb.Worksheets.Item('West Loop')
$ws.Range("J2").Value = 24
$ws.Range("J7").Value = 103

$ws = $wb.Worksheets.Item('East Side')
$ws.Range("J2").Value = 29
$ws.Range("J7").Value = 65

$ws = $wb.Worksheets.Item('Lincoln Square')
$ws.Range("J4").Value = 12
$ws.Range("J7").Value = 64

$ws = $wb.Worksheets.Item('North Center')
$ws.Range("J2").Value = 4
$ws.Range("J7").Value = 27

$ws = $wb.Worksheets.Item('West Town')
$ws.Range("J3").Value = 10
$ws.Range("J7").Value = 68

$ws = $wb.Worksheets.Item('Austin')
$ws.Range("J3").Value = 231
$ws.Range("J4").Value = 38
$ws.Range("J6").Value = 221
$ws.Range("J7").Value = 732

$ws = $wb.Worksheets.Item('Edgewater')
$ws.Range("J6").Value = 25
$ws.Range("J7").Value = 69

$ws = $wb.Worksheets.Item('Little Italy, UIC')
$ws.Range("J6").Value = 48
$ws.Range("J7").Value = 154

$ws = $wb.Worksheets.Item('Mckinley Park')
$ws.Range("J2").Value = 13
$ws.Range("J7").Value = 52

$ws = $wb.Worksheets.Item('Hyde Park')
$ws.Range("J4").Value = 11
$ws.Range("J7").Value = 105

$ws = $wb.Worksheets.Item('Logan Square')
$ws.Range("J2").Value = 24
$ws.Range("J7").Value = 111

$ws = $wb.Worksheets.Item('Clearing')
$ws.Range("J3").Value = 6
$ws.Range("J7").Value = 25

$ws = $wb.Worksheets.Item('Oakland')
$ws.Range("J3").Value = 11
$ws.Range("J7").Value = 42

$ws = $wb.Worksheets.Item('Old Town')
$ws.Range("J2").Value = 16
$ws.Range("J7").Value = 46

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range("J2").Value = 116
$ws.Range("J3").Value = 104
$ws.Range("J7").Value = 346
